$d = $word.ActiveDocument

# Locate the very end of the document body (before the final paragraph mark)
$endPos = $d.Content.End - 1
$r = $d.Range($endPos, $endPos)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:br w:type="page"/></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Problemas comunes encontrados y como solucionarlos (Windows + Git Bash)</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>1) Error: detected dubious ownership in repository (safe.directory)</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b w:val="0"/><w:sz w:val="22"/></w:rPr><w:t>Causa: Git detecta que la carpeta del repositorio pertenece a otro usuario/SID de Windows distinto al usuario actual. Por seguridad, bloquea operaciones hasta que marques la ruta como segura o ajustes permisos.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b w:val="0"/><w:sz w:val="22"/></w:rPr><w:t>Solucion rapida (recomendada): marcar la carpeta como segura:</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="120" w:after="120"/><w:ind w:left="360"/><w:shd w:val="clear" w:color="auto" w:fill="F2F2F2"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:eastAsia="Courier New"/><w:sz w:val="20"/></w:rPr><w:t>git config --global --add safe.directory "E:/Desarrollo de sofware 1/Desarrollo_Sofware1_2026"</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b w:val="0"/><w:sz w:val="22"/></w:rPr><w:t>Luego reintenta el comando que estabas ejecutando:</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="120" w:after="120"/><w:ind w:left="360"/><w:shd w:val="clear" w:color="auto" w:fill="F2F2F2"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:eastAsia="Courier New"/><w:sz w:val="20"/></w:rPr><w:t>git remote -v</w:t><w:br/><w:t># si hace falta:</w:t><w:br/><w:t># git remote add origin &lt;URL&gt;</w:t><w:br/><w:t># o:</w:t><w:br/><w:t># git remote set-url origin &lt;URL&gt;</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b w:val="0"/><w:sz w:val="22"/></w:rPr><w:t>Solucion alternativa: cambiar el propietario/permisos de la carpeta en Windows (evita que vuelva a pasar).</w:t></w:r></w:p><w:p/><w:p><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>2) Error 403: Permission denied (estas autenticado con otra cuenta)</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b w:val="0"/><w:sz w:val="22"/></w:rPr><w:t>Ejemplo: Permission to &lt;owner&gt;/&lt;repo&gt;.git denied to &lt;usuario&gt;. Causa: estas logueado en GitHub con una cuenta que no tiene permisos de escritura en ese repositorio, o el repositorio no existe/esta privado sin acceso.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b w:val="0"/><w:sz w:val="22"/></w:rPr><w:t>Soluciones:</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b w:val="0"/><w:sz w:val="22"/></w:rPr><w:t>A) Que el propietario te invite como collaborator (o te de acceso).</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b w:val="0"/><w:sz w:val="22"/></w:rPr><w:t>B) Cambiar el remote origin para apuntar a un repositorio en tu propia cuenta (donde si tengas permisos).</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b w:val="0"/><w:sz w:val="22"/></w:rPr><w:t>Ejemplo para cambiar el remote:</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="120" w:after="120"/><w:ind w:left="360"/><w:shd w:val="clear" w:color="auto" w:fill="F2F2F2"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:eastAsia="Courier New"/><w:sz w:val="20"/></w:rPr><w:t>git remote set-url origin https://github.com/TU_USUARIO/TU_REPO.git</w:t><w:br/><w:t>git remote -v</w:t><w:br/><w:t>git push -u origin main</w:t></w:r></w:p><w:p/><w:p><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>3) Cambiar de usuario en Git (HTTPS / Windows Credential Manager)</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b w:val="0"/><w:sz w:val="22"/></w:rPr><w:t>En Windows, Git puede reutilizar credenciales guardadas. Para forzar a iniciar sesion con otra cuenta, borra la credencial guardada para github.com.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b w:val="0"/><w:sz w:val="22"/></w:rPr><w:t>Borrar credencial guardada (Git Credential Manager):</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="120" w:after="120"/><w:ind w:left="360"/><w:shd w:val="clear" w:color="auto" w:fill="F2F2F2"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:eastAsia="Courier New"/><w:sz w:val="20"/></w:rPr><w:t>git credential-manager erase</w:t><w:br/><w:t>protocol=https</w:t><w:br/><w:t>host=github.com</w:t><w:br/></w:r></w:p><w:p><w:r><w:rPr><w:b w:val="0"/><w:sz w:val="22"/></w:rPr><w:t>Despues, al hacer push, Git te pedira credenciales nuevamente. Nota: GitHub no acepta contrasena, usa un PAT (token).</w:t></w:r></w:p><w:p/><w:p><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>4) Cambiar de usuario con GitHub CLI (gh) (si lo tienes instalado)</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="120" w:after="120"/><w:ind w:left="360"/><w:shd w:val="clear" w:color="auto" w:fill="F2F2F2"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:eastAsia="Courier New"/><w:sz w:val="20"/></w:rPr><w:t>gh auth status</w:t><w:br/><w:t>gh auth logout</w:t><w:br/><w:t>gh auth login</w:t><w:br/><w:t># luego:</w:t><w:br/><w:t># git push -u origin main</w:t></w:r></w:p><w:p/><w:p><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>5) Recomendado: usar SSH (evita tokens y cambios de usuario por HTTPS)</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b w:val="0"/><w:sz w:val="22"/></w:rPr><w:t>Generar llave SSH y agregarla a GitHub (en la cuenta correcta):</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="120" w:after="120"/><w:ind w:left="360"/><w:shd w:val="clear" w:color="auto" w:fill="F2F2F2"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:eastAsia="Courier New"/><w:sz w:val="20"/></w:rPr><w:t>ssh-keygen -t ed25519 -C "tu_correo@dominio.com"</w:t><w:br/><w:t>cat ~/.ssh/id_ed25519.pub</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b w:val="0"/><w:sz w:val="22"/></w:rPr><w:t>Cambiar el remote a SSH y probar:</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="120" w:after="120"/><w:ind w:left="360"/><w:shd w:val="clear" w:color="auto" w:fill="F2F2F2"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:eastAsia="Courier New"/><w:sz w:val="20"/></w:rPr><w:t>git remote set-url origin git@github.com:OWNER/REPO.git</w:t><w:br/><w:t>ssh -T git@github.com</w:t><w:br/><w:t>git push -u origin main</w:t></w:r></w:p><w:p/><w:p><w:r><w:rPr><w:b/><w:sz w:val="22"/></w:rPr><w:t>6) Warnings CRLF/LF (no bloquean el push)</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b w:val="0"/><w:sz w:val="22"/></w:rPr><w:t>Es normal al trabajar en Windows. Si quieres que Git lo gestione automaticamente:</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="120" w:after="120"/><w:ind w:left="360"/><w:shd w:val="clear" w:color="auto" w:fill="F2F2F2"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:eastAsia="Courier New"/><w:sz w:val="20"/></w:rPr><w:t>git config --global core.autocrlf true</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($xml)
